$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update displayed cell text/values (rows 2 and 3) ---
# Row 2: Server URL (hyperlink text), Project Name, PAT
$ws.Range("A2").Value = "http://128.0.0.1/TestCollection"
$ws.Range("B2").Value = "project1"
$ws.Range("C2").Value = "adad87adad8ds4449m434344mmnbnbb43434"

# Row 3: Server URL (hyperlink text), Project Name, PAT
$ws.Range("A3").Value = "http://128.0.0.1/TestCollection"
$ws.Range("B3").Value = "project2"
$ws.Range("C3").Value = "adad87adad8ds4449m434344mmnbnbb43434"

# --- Remove the now-unneeded extra rows (old devserver/qaserver rows 4 & 5) ---
$ws.Rows("4:5").Delete()

# --- Tidy up the hyperlinks collection: the engine leaves stale hyperlink
# entries (still addressed at the now-deleted A4/A5) after the row delete
# above, so rebuild the two that should remain (A2, A3) pointing at the
# same external target they always pointed to. ---
$hyperlinkTarget = "http://172.191.4.85/TestCollection"
$savedStyleA2 = $ws.Range("A2").Style
$savedStyleA3 = $ws.Range("A3").Style

$ws.Range("A2").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), $hyperlinkTarget)
$ws.Hyperlinks.Add($ws.Range("A3"), $hyperlinkTarget)

# Adding a hyperlink re-applies the built-in "Hyperlink" cell style; restore
# the original style reference used throughout the sheet.
$ws.Range("A2").Style = $savedStyleA2
$ws.Range("A3").Style = $savedStyleA3

# Reset selection to the top of the sheet now that the old selection (B5)
# no longer exists.
$ws.Range("A1").Select()
